$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "preprocessing session 4 and trial by trial pipeline session 2"
$ws.Range("B14").Value = "6, 45 + 5, 15 "
